$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = 1995
$ws.Range("B2").Value = 23.98030410786858
$ws.Range("A3").Value = 1996
$ws.Range("B3").Value = 23.56255995421238
$ws.Range("A4").Value = 1997
$ws.Range("B4").Value = 23.12549774171843
$ws.Range("A5").Value = 1998
$ws.Range("B5").Value = 22.69324139491783
$ws.Range("A6").Value = 1999
$ws.Range("B6").Value = 22.25307221153614
$ws.Range("A7").Value = 2000
$ws.Range("B7").Value = 21.89644669448948
$ws.Range("A8").Value = 2001
$ws.Range("B8").Value = 21.52499277652478
$ws.Range("A9").Value = 2002
$ws.Range("B9").Value = 21.15599090857465
$ws.Range("A10").Value = 2003
$ws.Range("B10").Value = 20.84001722660343
$ws.Range("A11").Value = 2004
$ws.Range("B11").Value = 20.8420513580753
$ws.Range("A12").Value = 2005
$ws.Range("B12").Value = 20.97507412129914
$ws.Range("A13").Value = 2006
$ws.Range("B13").Value = 20.90415978963375
$ws.Range("A14").Value = 2007
$ws.Range("B14").Value = 21.02701795652256
$ws.Range("A15").Value = 2008
$ws.Range("B15").Value = 21.39020486169563
$ws.Range("A16").Value = 2009
$ws.Range("B16").Value = 20.05934924949733
$ws.Range("A17").Value = 2010
$ws.Range("B17").Value = 20.96253299247605
$ws.Range("A18").Value = 2011
$ws.Range("B18").Value = 20.25246513123966
$ws.Range("A19").Value = 2012
$ws.Range("B19").Value = 19.49103920938649
$ws.Range("A20").Value = 2013
$ws.Range("B20").Value = 19.88881356008027
$ws.Range("A21").Value = 2014
$ws.Range("B21").Value = 19.74535937402517
$ws.Range("A22").Value = 2015
$ws.Range("B22").Value = 19.87884121989851
$ws.Range("A23").Value = 2016
$ws.Range("B23").Value = 19.95872818965945

$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = 2017
$ws.Range("B2").Value = 19.60072106797074
$ws.Range("A3").Value = 2018
$ws.Range("B3").Value = 19.15464086645493
$ws.Range("A4").Value = 2019
$ws.Range("B4").Value = 18.92535626220639
$ws.Range("A5").Value = 2020
$ws.Range("B5").Value = 18.7619460158259
$ws.Range("A6").Value = 2021
$ws.Range("B6").Value = 18.62508240777986

$ws = $wb.Worksheets.Item(3)
$ws.Rows.Item(29).Delete()
$ws.Range("A2").Value = 1995
$ws.Range("B2").Value = 23.60359831301482
$ws.Range("A3").Value = 1996
$ws.Range("B3").Value = 23.34717677364178
$ws.Range("A4").Value = 1997
$ws.Range("B4").Value = 23.04580400479329
$ws.Range("A5").Value = 1998
$ws.Range("B5").Value = 22.76398591110084
$ws.Range("A6").Value = 1999
$ws.Range("B6").Value = 22.36357599041737
$ws.Range("A7").Value = 2000
$ws.Range("B7").Value = 22.03886908665579
$ws.Range("A8").Value = 2001
$ws.Range("B8").Value = 21.80712655575084
$ws.Range("A9").Value = 2002
$ws.Range("B9").Value = 21.27252028004369
$ws.Range("A10").Value = 2003
$ws.Range("B10").Value = 20.96139389822509
$ws.Range("A11").Value = 2004
$ws.Range("B11").Value = 20.95260205134854
$ws.Range("A12").Value = 2005
$ws.Range("B12").Value = 21.04287426514862
$ws.Range("A13").Value = 2006
$ws.Range("B13").Value = 21.03817319395989
$ws.Range("A14").Value = 2007
$ws.Range("B14").Value = 21.06668981579421
$ws.Range("A15").Value = 2008
$ws.Range("B15").Value = 21.2095482179747
$ws.Range("A16").Value = 2009
$ws.Range("B16").Value = 20.51393243399232
$ws.Range("A17").Value = 2010
$ws.Range("B17").Value = 20.15952077884339
$ws.Range("A18").Value = 2011
$ws.Range("B18").Value = 20.29021369707276
$ws.Range("A19").Value = 2012
$ws.Range("B19").Value = 19.69991137232896
$ws.Range("A20").Value = 2013
$ws.Range("B20").Value = 19.90448812452366
$ws.Range("A21").Value = 2014
$ws.Range("B21").Value = 19.7726284689958
$ws.Range("A22").Value = 2015
$ws.Range("B22").Value = 19.82329899438644
$ws.Range("A23").Value = 2016
$ws.Range("B23").Value = 19.92649517983924
$ws.Range("A24").Value = 2017
$ws.Range("B24").Value = 19.71050122427598
$ws.Range("A25").Value = 2018
$ws.Range("B25").Value = 19.95694092041195
$ws.Range("A26").Value = 2019
$ws.Range("B26").Value = 20.00884064375798
$ws.Range("A27").Value = 2020
$ws.Range("B27").Value = 19.80605197471866
$ws.Range("A28").Value = 2021
$ws.Range("B28").Value = 20.3170655497136

$ws = $wb.Worksheets.Item(4)
$ws.Range("A2").Value = 2022
$ws.Range("B2").Value = 19.09816268486769
$ws.Range("A3").Value = 2023
$ws.Range("B3").Value = 19.20846279063443
$ws.Range("A4").Value = 2024
$ws.Range("B4").Value = 19.23260074467692
$ws.Range("A5").Value = 2025
$ws.Range("B5").Value = 19.23610997346741
$ws.Range("A6").Value = 2026
$ws.Range("B6").Value = 19.21901970489343

Write-Host "All updates applied."
